$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "IK106" -- refresh the "gewijzigd" (last-modified) timestamp column
# for the existing 5 data rows (rows 2-6), column E.
# ---------------------------------------------------------------------------
$wsIK106 = $wb.Worksheets.Item("IK106")
foreach ($r in 2..6) {
    $wsIK106.Cells.Item($r, 5).Value = 45693.44442743056
}

# ---------------------------------------------------------------------------
# Sheet "IK96" -- refresh the "gewijzigd" timestamp column for the existing
# 3 data rows (rows 2-4), column E.
# ---------------------------------------------------------------------------
$wsIK96 = $wb.Worksheets.Item("IK96")
foreach ($r in 2..4) {
    $wsIK96.Cells.Item($r, 5).Value = 45699.61820989798
}

# ---------------------------------------------------------------------------
# Sheet "P100" -- model coefficients were recomputed with a new data point
# added (a new row 5), the B/C values on rows 2-4 shifted down one data
# point, and every row's "gewijzigd" timestamp was refreshed.
# ---------------------------------------------------------------------------
$wsP100 = $wb.Worksheets.Item("P100")

# Row 2: offset recomputed; B2/D2 unchanged.
$wsP100.Cells.Item(2, 3).Value = [double]"-3.011893641386587e-114"
$wsP100.Cells.Item(2, 5).Value = 45699.44560905093

# Row 3: now carries what used to be an intermediate point.
$wsP100.Cells.Item(3, 2).Value = 43089
$wsP100.Cells.Item(3, 3).Value = -0.002418679012150875
$wsP100.Cells.Item(3, 5).Value = 45699.44560905093

# Row 4: shifted down from the old row 3 values.
$wsP100.Cells.Item(4, 2).Value = 44237
$wsP100.Cells.Item(4, 3).Value = -0.1151145706160745
$wsP100.Cells.Item(4, 5).Value = 45699.44560905093

# New row 5: carries what used to be the old row 4 (B/date), with the
# offset recomputed, matching the formatting used by the rows above it.
$wsP100.Cells.Item(5, 1).Value = 3
$wsP100.Cells.Item(5, 2).Value = 45203
$wsP100.Cells.Item(5, 3).Value = -0.1442102341436769
$wsP100.Cells.Item(5, 4).Value = -0.0003
$wsP100.Cells.Item(5, 5).Value = 45699.44560905093

# Match the look of column A's index cells (bold, boxed, centered) ...
$a5 = $wsP100.Cells.Item(5, 1)
$a5.BorderAround(1)
$a5.Font.Bold = $true
$a5.HorizontalAlignment = -4108
$a5.VerticalAlignment = -4160

# ... and the date-time number format used by columns B and E elsewhere.
$wsP100.Cells.Item(5, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsP100.Cells.Item(5, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
